$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to remain text so values like "0.635" or "1.00"
# are not re-interpreted by Excel as numbers (losing trailing zeros / formatting).
$priceRange = $ws.Range("D2:D51")
$priceRange.Style = "Normal"
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "62.007.70"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "3.427.07"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "410.69"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").Value = "130.12"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").Value = "0.635"
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "0.737"
$ws.Range("E9").Value = "  -2.57%  "

# Row 10
$ws.Range("E10").Value = "  -1.91%  "

# Row 11
$ws.Range("D11").Value = "43.70"
$ws.Range("E11").Value = "  +1.19%  "

# Row 12
$ws.Range("E12").Value = "  +15.87%  "

# Row 13
$ws.Range("D13").Value = "9.37"
$ws.Range("E13").Value = "  +5.26%  "

# Row 14
$ws.Range("D14").Value = "3.969.98"
$ws.Range("E14").Value = "  -0.22%  "

# Row 15
$ws.Range("E15").Value = "  +0.29%  "

# Row 16
$ws.Range("D16").Value = "21.27"
$ws.Range("E16").Value = "  +3.70%  "

# Row 17
$ws.Range("D17").Value = "3.418.68"
$ws.Range("E17").Value = "  -0.87%  "

# Row 18
$ws.Range("D18").Value = "12.39"
$ws.Range("E18").Value = "  +7.07%  "

# Row 19
$ws.Range("E19").Value = "  +2.89%  "

# Row 20
$ws.Range("D20").Value = "61.920.26"
$ws.Range("E20").Value = "  -0.35%  "

# Row 21
$ws.Range("D21").Value = "513.58"
$ws.Range("E21").Value = "  +31.74%  "

# Row 22
$ws.Range("D22").Value = "92.52"
$ws.Range("E22").Value = "  +4.20%  "

# Row 23
$ws.Range("E23").Value = "  +4.10%  "

# Row 24
$ws.Range("D24").Value = "13.48"
$ws.Range("E24").Value = "  +0.76%  "

# Row 25
$ws.Range("D25").Value = "3.33"
$ws.Range("E25").Value = "  +3.54%  "

# Row 26
$ws.Range("D26").Value = "34.87"
$ws.Range("E26").Value = "  +8.75%  "

# Row 27
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  +9.63%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  -0.21%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "12.18"
$ws.Range("E29").Value = "  +3.16%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").Value = "  -2.65%  "

# Row 31
$ws.Range("E31").Value = "  -1.42%  "

# Row 32
$ws.Range("E32").Value = "  -2.00%  "

# Row 33
$ws.Range("D33").Value = "42.03"
$ws.Range("E33").Value = "  -4.75%  "

# Row 34
$ws.Range("D34").Value = "59.21"
$ws.Range("E34").Value = "  +13.02%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").Value = "0.0501"
$ws.Range("E36").Value = "  +1.73%  "

# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "0.140"
$ws.Range("E37").Value = "  +5.36%  "

# Row 38
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39
$ws.Range("D39").Value = "3.47"
$ws.Range("E39").Value = "  +2.53%  "

# Row 40
$ws.Range("E40").Value = "  +18.15%  "

# Row 41
$ws.Range("D41").Value = "147.93"
$ws.Range("E41").Value = "  +4.80%  "

# Row 42
$ws.Range("E42").Value = "  +7.49%  "

# Row 43
$ws.Range("E43").Value = "  +0.65%  "

# Row 44
$ws.Range("E44").Value = "  +1.89%  "

# Row 45
$ws.Range("E45").Value = "  +7.99%  "

# Row 46
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  +22.55%  "

# Row 47
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "16.75"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "121.28"
$ws.Range("E48").Value = "  +27.83%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "23.14"
$ws.Range("E49").Value = "  +4.45%  "

# Row 50
$ws.Range("D50").Value = "0.148"
$ws.Range("E50").Value = "  +20.23%  "

# Row 51
$ws.Range("D51").Value = "2.141.85"
$ws.Range("E51").Value = "  +0.81%  "

# Restore default style on the Price column (values already committed as text).
$priceRange.Style = "Normal"
